$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (D1 direction) value swaps ---
# Rows 2-5: S -> W
$ws.Range("C2").Value = "W"
$ws.Range("C3").Value = "W"
$ws.Range("C4").Value = "W"
$ws.Range("C5").Value = "W"

# Rows 6-9: W -> S
$ws.Range("C6").Value = "S"
$ws.Range("C7").Value = "S"
$ws.Range("C8").Value = "S"
$ws.Range("C9").Value = "S"

# Rows 10-13: S -> W
$ws.Range("C10").Value = "W"
$ws.Range("C11").Value = "W"
$ws.Range("C12").Value = "W"
$ws.Range("C13").Value = "W"

# Rows 54-57: S -> W
$ws.Range("C54").Value = "W"
$ws.Range("C55").Value = "W"
$ws.Range("C56").Value = "W"
$ws.Range("C57").Value = "W"

# Rows 86-91: S -> W
$ws.Range("C86").Value = "W"
$ws.Range("C87").Value = "W"
$ws.Range("C88").Value = "W"
$ws.Range("C89").Value = "W"
$ws.Range("C90").Value = "W"
$ws.Range("C91").Value = "W"

# --- Apply a sort (Date ascending) over A2:G96 so the worksheet records a
# sortState / sortCondition, matching the re-sorted view used for Figure 3 ---
$sortRange = $ws.Range("A2:G96")
$keyRange = $ws.Range("A2:A96")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# --- View state: active cell moved to C9, scrolled so row 13 is at the top ---
$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 13
